# "Se Agregan adicionales de tablas"
# Populates the "Nueva Tabla" column (F) on Hoja2 for the newly documented
# tables: PLFOM3, PLFOM5, PLLIM1-4, PLODP1-5 and PLRCI1. Each block of rows
# belonging to one TABLA code gets the matching T_PRD_DET_* label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

$ws.Range("F1264:F1270").Value = "T_PRD_DET_FOM3_DetalleCostosIndirectosFabricacionFormulasMaestras"
$ws.Range("F1271:F1276").Value = "T_PRD_DET_FOM5_DetalleAnalisis"
$ws.Range("F1277:F1288").Value = "T_PRD_DET_LIM1_DetalleListadoFormulaMaestra"
$ws.Range("F1289:F1299").Value = "T_PRD_DET_LIM2_DetalleListadoMateriales"
$ws.Range("F1300:F1306").Value = "T_PRD_DET_LIM3_DetalleListadoMaterialesOperaciones"
$ws.Range("F1307:F1318").Value = "T_PRD_DET_LIM4_DetalleListadoMaterialesCIF"
$ws.Range("F1319:F1337").Value = "T_PRD_DET_ODP1_DetalleOrdenProduccion"
$ws.Range("F1338:F1354").Value = "T_PRD_DET_ODP2_DetalleOrdenOperacion"
$ws.Range("F1355:F1370").Value = "T_PRD_DET_ODP3_DetalleCIF"
$ws.Range("F1371:F1406").Value = "T_PRD_DET_ODP4_ResumenOP"
$ws.Range("F1407:F1414").Value = "T_PRD_DET_ODP5_DetalleAnalisis"
$ws.Range("F1415:F1444").Value = "T_PRD_DET_RCI1_DetalleRegistroCIF"

# Restore the view state the author had when saving (best effort - cosmetic).
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("F1494").Select()
